$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Two new match rows are appended at the bottom of the data table (rows 104
# and 105), matching the structure of every existing row (Indice, pais,
# torneio, temporada, data_partida, home, home_ft_gols, away, away_ft_gols,
# odds columns, url_partida).
# ---------------------------------------------------------------------------

# Row 104: Ameliano vs Guarani
$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = "paraguay"
$ws.Cells.Item(104, 3).Value = "primera-division"
$ws.Cells.Item(104, 4).Value = 2023
$ws.Cells.Item(104, 5).Value = 45231.91666666666
$ws.Cells.Item(104, 6).Value = "Ameliano"
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = "Guarani"
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 2.15
$ws.Cells.Item(104, 11).Value = "29/10/2023 09:04"
$ws.Cells.Item(104, 12).Value = 2.23
$ws.Cells.Item(104, 13).Value = "01/11/2023 21:50"
$ws.Cells.Item(104, 14).Value = 3.52
$ws.Cells.Item(104, 15).Value = "29/10/2023 09:04"
$ws.Cells.Item(104, 16).Value = 3.48
$ws.Cells.Item(104, 17).Value = "01/11/2023 21:50"
$ws.Cells.Item(104, 18).Value = 3.38
$ws.Cells.Item(104, 19).Value = "29/10/2023 09:04"
$ws.Cells.Item(104, 20).Value = 3.36
$ws.Cells.Item(104, 21).Value = "01/11/2023 21:50"
$ws.Cells.Item(104, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/sportivo-ameliano-guarani/YFqsnGDD/"

# Row 105: General Caballero JLM vs Nacional Asuncion
$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = "paraguay"
$ws.Cells.Item(105, 3).Value = "primera-division"
$ws.Cells.Item(105, 4).Value = 2023
$ws.Cells.Item(105, 5).Value = 45232.02083333334
$ws.Cells.Item(105, 6).Value = "General Caballero JLM"
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = "Nacional Asuncion"
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 4
$ws.Cells.Item(105, 11).Value = "29/10/2023 09:04"
$ws.Cells.Item(105, 12).Value = 4.12
$ws.Cells.Item(105, 13).Value = "02/11/2023 00:21"
$ws.Cells.Item(105, 14).Value = 3.55
$ws.Cells.Item(105, 15).Value = "29/10/2023 09:04"
$ws.Cells.Item(105, 16).Value = 3.34
$ws.Cells.Item(105, 17).Value = "02/11/2023 00:21"
$ws.Cells.Item(105, 18).Value = 1.94
$ws.Cells.Item(105, 19).Value = "29/10/2023 09:04"
$ws.Cells.Item(105, 20).Value = 2.03
$ws.Cells.Item(105, 21).Value = "02/11/2023 00:21"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/general-caballero-jlm-nacional-asuncion/zLNPAajE/"

# The "temporada" column holds the season as TEXT ("2023"), even though it
# looks numeric, matching every other row in the column. Marking the cells
# as Text first means the numeric-looking string is kept as a real string
# (not auto-converted to a number); clearing the format afterwards drops the
# now-unneeded explicit "Text" number format again so the cell is left on
# the default style, same as the rest of the column.
$ws.Range("D104:D105").NumberFormat = "@"
$ws.Cells.Item(104, 4).Value = "2023"
$ws.Cells.Item(105, 4).Value = "2023"
$ws.Range("D104:D105").ClearFormats()

# Re-apply formatting to match the rest of the table: column A (Indice) uses
# the bold/bordered/centered header-like style, column E (data_partida) uses
# the custom date-time number format. Copy the formatting straight from the
# row above (row 103) so the exact same style is reused instead of a new one
# being minted.
$ws.Range("A103").Copy()
$ws.Range("A104:A105").PasteSpecial(-4122)

$ws.Range("E103").Copy()
$ws.Range("E104:E105").PasteSpecial(-4122)

$excel.CutCopyMode = 0
